$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 269-274
$data = @(
    @(44192, 2232, 76, 431, 15, 234, 4),
    @(44193, 2259, 27, 431, 0, 225, -9),
    @(44194, 2257, -2, 433, 2, 231, 6),
    @(44195, 2271, 14, 417, -16, 240, 9),
    @(44196, 2323, 52, 429, 12, 258, 18),
    @(44197, 2280, -43, 412, -17, 246, -12)
)

$startRow = 269
$endRow = 274

# Carry column A's existing date style (s="1") down onto the new rows
# instead of using NumberFormat (which would mint a brand-new style entry).
$ws.Range("A268").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($r -eq 269) {
        # Matches the source workbook's off-by-one formula (excludes B269 itself)
        $ws.Cells.Item($r, 4).Formula = "=AVERAGE(B263:B268)"
    } else {
        $ws.Cells.Item($r, 4).Formula = "=AVERAGE(B" + ($r - 6) + ":B" + $r + ")"
    }
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# Update view: pane and selection
$ws.Application.ActiveWindow.SelectedSheets.Item(1).Activate()
$ws.Range("G276:H277").Select()
